$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column K, copying the header style from A1
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "intervention_type"

# Add data for the new column
$ws.Range("K2").Value = "OTHER"
$ws.Range("K3").Value = "OTHER"
$ws.Range("K4").Value = "DRUG"
